$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename column header, add new "EventDate" header, bold @ size 11 ---
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Size = 11
$ws.Range("A1").Value = "FixtureName"
$ws.Range("B1").Value = "EventDate"

# --- Data block: fixture name (col A) + event date serials (col B), with 3 new fixtures inserted ---
$data = New-Object "object[,]" 45,2
$data[0,0] = 'Robbie Williams Live 2025 - Saturday'
$data[0,1] = 45815.6875
$data[1,0] = 'Robbie Williams Live 2025 (Friday)'
$data[1,1] = 45814.6875
$data[2,0] = 'Arsenal v Newcastle United'
$data[2,1] = 45795.58333333334
$data[3,0] = 'Arsenal Women v Manchester United Women'
$data[3,1] = 45787.47916666666
$data[4,0] = 'Arsenal v A.F.C. Bournemouth'
$data[4,1] = 45780.6875
$data[5,0] = 'Arsenal v Crystal Palace'
$data[5,1] = 45773.6875
$data[6,0] = 'Arsenal v Brentford'
$data[6,1] = 45759.6875
$data[7,0] = 'Arsenal v Real Madrid'
$data[7,1] = 45755.79166666666
$data[8,0] = 'Arsenal v Fulham'
$data[8,1] = 45748.78125
$data[9,0] = 'Arsenal Women v Real Madrid Women'
$data[9,1] = 45742.83333333334
$data[10,0] = 'Arsenal Women v Liverpool Women'
$data[10,1] = 45738.72916666666
$data[11,0] = 'Arsenal v Chelsea'
$data[11,1] = 45732.5625
$data[12,0] = 'Arsenal v PSV'
$data[12,1] = 45728.83333333334
$data[13,0] = 'Arsenal v West Ham United'
$data[13,1] = 45710.625
$data[14,0] = 'Arsenal Women v Tottenham Hotspur Women'
$data[14,1] = 45704.52083333334
$data[15,0] = 'Arsenal v Manchester City'
$data[15,1] = 45690.6875
$data[16,0] = 'Arsenal v Dinamo Zagreb'
$data[16,1] = 45679.83333333334
$data[17,0] = 'Arsenal v Aston Villa'
$data[17,1] = 45675.72916666666
$data[18,0] = 'Arsenal v Tottenham Hotspur'
$data[18,1] = 45672.83333333334
$data[19,0] = 'Arsenal v Manchester United'
$data[19,1] = 45669.625
$data[20,0] = 'Arsenal v Newcastle United'
$data[20,1] = 45664.83333333334
$data[21,0] = 'Arsenal v Ipswich Town'
$data[21,1] = 45653.84375
$data[22,0] = 'Arsenal v Crystal Palace'
$data[22,1] = 45644.8125
$data[23,0] = 'Arsenal Women v FC Bayern Munich Women'
$data[23,1] = 45644.83333333334
$data[24,0] = 'Arsenal v Everton'
$data[24,1] = 45640.625
$data[25,0] = 'Arsenal v AS Monaco'
$data[25,1] = 45637.83333333334
$data[26,0] = 'Arsenal Women v Aston Villa Women'
$data[26,1] = 45634.58333333334
$data[27,0] = 'Arsenal v Manchester United'
$data[27,1] = 45630.84375
$data[28,0] = 'Arsenal v Nottingham Forest'
$data[28,1] = 45619.625
$data[29,0] = 'Arsenal Women v Juventus Women'
$data[29,1] = 45617.83333333334
$data[30,0] = 'Arsenal Women v Brighton Women'
$data[30,1] = 45604.79166666666
$data[31,0] = 'Arsenal v Liverpool'
$data[31,1] = 45592.6875
$data[32,0] = 'Arsenal v Shakhtar Donetsk'
$data[32,1] = 45587.79166666666
$data[33,0] = 'Arsenal Women v Vålerenga Women'
$data[33,1] = 45581.79166666666
$data[34,0] = 'Arsenal Women v Chelsea Women'
$data[34,1] = 45577.53125
$data[35,0] = 'Arsenal Women v Everton Women'
$data[35,1] = 45571.54166666666
$data[36,0] = 'Arsenal v Southampton'
$data[36,1] = 45570.58333333334
$data[37,0] = 'Arsenal v Paris Saint-Germain'
$data[37,1] = 45566.79166666666
$data[38,0] = 'Arsenal v Leicester City'
$data[38,1] = 45563.58333333334
$data[39,0] = 'Arsenal v Bolton Wanderers'
$data[39,1] = 45560.78125
$data[40,0] = 'Arsenal Women v Manchester City Women'
$data[40,1] = 45557.47916666666
$data[41,0] = 'Arsenal v Brighton'
$data[41,1] = 45535.47916666666
$data[42,0] = 'Arsenal v Wolves'
$data[42,1] = 45521.58333333334
$data[43,0] = 'Arsenal v Olympique Lyonnais'
$data[43,1] = 45515.54166666666
$data[44,0] = 'Arsenal v Bayer 04 Leverkusen'
$data[44,1] = 45511.70833333334
$ws.Range("A2:B46").Value = $data

# --- Date/time number format on the new EventDate column ---
$ws.Range("B2:B46").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# --- Column B width (column A already has the correct width from before.xlsx) ---
$ws.Columns.Item(2).ColumnWidth = 19.417

# --- Selection, matching the saved view state ---
$ws.Range("A1:B12").Select()
